{"js": "// Update the Visit-level N (%) summary counts/percentages in the table\n// (Q4 2021 refresh): visits 1-5, both the \"formatted_stats.x\" and\n// \"formatted_stats.y\" columns.\nconst replacements = [\n  [\"307 (100.0)\", \"329 (100.0)\"],\n  [\"281 (100.0)\", \"303 (100.0)\"], // appears twice (visit 1 & visit 2, y column)\n  [\"281 (91.8)\", \"303 (92.4)\"],\n  [\"205 (69.5)\", \"225 (69.4)\"],\n  [\"205 (75.4)\", \"225 (75.3)\"],\n  [\"157 (55.9)\", \"169 (54.9)\"],\n  [\"157 (60.2)\", \"169 (59.7)\"],\n  [\"143 (54.0)\", \"152 (53.5)\"],\n  [\"143 (57.9)\", \"152 (57.6)\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the Visit-level N (%) summary counts/percentages in the table\n# (Q4 2021 refresh): visits 1-5, both the \"formatted_stats.x\" and\n# \"formatted_stats.y\" columns.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"307 (100.0)\"; New = \"329 (100.0)\" },\n    @{ Old = \"281 (100.0)\"; New = \"303 (100.0)\" }, # appears twice (visit 1 & visit 2, y column)\n    @{ Old = \"281 (91.8)\";  New = \"303 (92.4)\" },\n    @{ Old = \"205 (69.5)\";  New = \"225 (69.4)\" },\n    @{ Old = \"205 (75.4)\";  New = \"225 (75.3)\" },\n    @{ Old = \"157 (55.9)\";  New = \"169 (54.9)\" },\n    @{ Old = \"157 (60.2)\";  New = \"169 (59.7)\" },\n    @{ Old = \"143 (54.0)\";  New = \"152 (53.5)\" },\n    @{ Old = \"143 (57.9)\";  New = \"152 (57.6)\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
